# Reorders data rows 2-34 on the active sheet according to $rowMap
# ($rowMap[newRow] = oldRow). Values for every data row (columns A-T) are
# first snapshotted, then written back in the new arrangement so reads
# never observe a partially-overwritten row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 34
$firstCol = 1
$lastCol = 20

$rowMap = @{
    2 = 20
    3 = 31
    4 = 32
    5 = 26
    6 = 5
    7 = 9
    8 = 19
    9 = 4
    10 = 3
    11 = 10
    12 = 11
    13 = 6
    14 = 7
    15 = 8
    16 = 16
    17 = 17
    18 = 23
    19 = 24
    20 = 21
    21 = 22
    22 = 2
    23 = 13
    24 = 27
    25 = 12
    26 = 33
    27 = 34
    28 = 18
    29 = 25
    30 = 28
    31 = 29
    32 = 30
    33 = 14
    34 = 15
}

# 1) Snapshot every data row (array of per-column Value2 reads).
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @()
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $rowVals += $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# 2) Write the snapshot back out in the permuted order.
for ($newRow = $firstRow; $newRow -le $lastRow; $newRow++) {
    $oldRow = $rowMap[$newRow]
    $rowVals = $snapshot[$oldRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($newRow, $c).Value2 = $rowVals[$c - 1]
    }
}
